$wb = $excel.ActiveWorkbook

# --- SA sheet: marks updates for rows 8 and 9 ---
$sa = $wb.Worksheets.Item("SA")
$sa.Activate()
$sa.Range("D8").Value = 40
$sa.Range("I8").Value = 8
$sa.Range("J8").Value = 7
$sa.Range("D9").Value = 30
$sa.Range("H9").Value = 10
$sa.Range("I9").Value = 10
$sa.Range("E14").Select()

# --- OS sheet: clear the assignment score that drove a failing grade ---
$os = $wb.Worksheets.Item("OS")
$os.Activate()
$os.Range("G9").ClearContents()
$os.Range("G11").Select()

# --- DCN sheet: record new marks for quiz/assignment/mid rows ---
$dcn = $wb.Worksheets.Item("DCN")
$dcn.Activate()
$dcn.Range("G7").Value = 18
$dcn.Range("G8").Value = 18
$dcn.Range("G10").Value = 85
$dcn.Range("G12").Select()

# DCN ends up as the active / selected tab in the saved workbook
$dcn.Activate()
